$d = $word.ActiveDocument

# Locate the "Appendix: Quick prototype" Heading2 paragraph and the
# "Appendix: Links" Heading2 paragraph that follows the prototype section,
# then remove everything in between (inclusive of the prototype heading,
# the "Figure: PDF page 1" text and the embedded image paragraph), while
# leaving the blank paragraph that precedes the prototype heading intact.

$startPara = $null
$endPara = $null

foreach ($p in $d.Paragraphs) {
    $txt = $p.Range.Text
    if ($startPara -eq $null) {
        if ($txt -like "Appendix: Quick prototype*") {
            $startPara = $p
        }
    } elseif ($endPara -eq $null) {
        if ($txt -like "Appendix: Links*") {
            $endPara = $p
            break
        }
    }
}

if ($startPara -ne $null -and $endPara -ne $null) {
    $rng = $d.Range($startPara.Range.Start, $endPara.Range.Start)
    $rng.Delete()
}
